$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.088846
$ws.Range("H2").Value = 30.266538
$ws.Range("I2").Value = 0.1151445838515654
$ws.Range("J2").Value = 0.1151445838515654
$ws.Range("M2").Value = 3.303267
$ws.Range("N2").Value = 9.909801000000002
$ws.Range("O2").Value = 0.03362563178859915
$ws.Range("P2").Value = 0.03362563178859915
$ws.Range("Q2").Value = 33.32615205988201
$ws.Range("R2").Value = 299.9353685389381
$ws.Range("S2").Value = 0.003871809379044219
$ws.Range("T2").Value = 0.003871809379044219
$ws.Range("G3").Value = 10.088846
$ws.Range("H3").Value = 30.266538
$ws.Range("I3").Value = 0.1151445838515654
$ws.Range("J3").Value = 0.1151445838515654
$ws.Range("M3").Value = 37.82684066666667
$ws.Range("O3").Value = 0.3850586149964086
$ws.Range("P3").Value = 0.3850586149964086
$ws.Range("Q3").Value = 381.6291701525374
$ws.Range("R3").Value = 3434.662531372836
$ws.Range("S3").Value = 0.04433741398222162
$ws.Range("T3").Value = 0.04433741398222162
$ws.Range("G4").Value = 10.088846
$ws.Range("H4").Value = 30.266538
$ws.Range("I4").Value = 0.1151445838515654
$ws.Range("J4").Value = 0.1151445838515654
$ws.Range("M4").Value = 9.149395999999999
$ws.Range("N4").Value = 27.448188
$ws.Range("O4").Value = 0.09313634682999644
$ws.Range("P4").Value = 0.09313634682999644
$ws.Range("Q4").Value = 92.306847237016
$ws.Range("R4").Value = 830.761625133144
$ws.Range("S4").Value = 0.010724145897195
$ws.Range("T4").Value = 0.010724145897195
$ws.Range("G5").Value = 10.088846
$ws.Range("H5").Value = 30.266538
$ws.Range("I5").Value = 0.1151445838515654
$ws.Range("J5").Value = 0.1151445838515654
$ws.Range("M5").Value = 47.95707433333333
$ws.Range("N5").Value = 143.871223
$ws.Range("O5").Value = 0.4881794063849957
$ws.Range("P5").Value = 0.4881794063849957
$ws.Range("Q5").Value = 483.8315375595527
$ws.Range("R5").Value = 4354.483838035973
$ws.Range("S5").Value = 0.05621121459310457
$ws.Range("T5").Value = 0.05621121459310457
$ws.Range("I6").Value = 0.4327250566572728
$ws.Range("J6").Value = 0.4327250566572729
$ws.Range("M6").Value = 3.303267
$ws.Range("N6").Value = 9.909801000000002
$ws.Range("O6").Value = 0.03362563178859915
$ws.Range("P6").Value = 0.03362563178859915
$ws.Range("Q6").Value = 125.243068808792
$ws.Range("R6").Value = 1127.187619279128
$ws.Range("S6").Value = 0.01455065342085816
$ws.Range("T6").Value = 0.01455065342085817
$ws.Range("I7").Value = 0.4327250566572728
$ws.Range("J7").Value = 0.4327250566572729
$ws.Range("M7").Value = 37.82684066666667
$ws.Range("O7").Value = 0.3850586149964086
$ws.Range("P7").Value = 0.3850586149964086
$ws.Range("S7").Value = 0.1666245109906919
$ws.Range("T7").Value = 0.1666245109906919
$ws.Range("I8").Value = 0.4327250566572728
$ws.Range("J8").Value = 0.4327250566572729
$ws.Range("M8").Value = 9.149395999999999
$ws.Range("N8").Value = 27.448188
$ws.Range("O8").Value = 0.09313634682999644
$ws.Range("P8").Value = 0.09313634682999644
$ws.Range("Q8").Value = 346.8985197947626
$ws.Range("R8").Value = 3122.086678152864
$ws.Range("S8").Value = 0.04030243095886162
$ws.Range("T8").Value = 0.04030243095886163
$ws.Range("I9").Value = 0.4327250566572728
$ws.Range("J9").Value = 0.4327250566572729
$ws.Range("M9").Value = 47.95707433333333
$ws.Range("N9").Value = 143.871223
$ws.Range("O9").Value = 0.4881794063849957
$ws.Range("P9").Value = 0.4881794063849957
$ws.Range("Q9").Value = 1818.28812524026
$ws.Range("R9").Value = 16364.59312716234
$ws.Range("S9").Value = 0.2112474612868611
$ws.Range("T9").Value = 0.2112474612868611
$ws.Range("G10").Value = 15.69885766666667
$ws.Range("H10").Value = 47.096573
$ws.Range("I10").Value = 0.1791719719949428
$ws.Range("J10").Value = 0.1791719719949428
$ws.Range("M10").Value = 3.303267
$ws.Range("N10").Value = 9.909801000000002
$ws.Range("O10").Value = 0.03362563178859915
$ws.Range("P10").Value = 0.03362563178859915
$ws.Range("Q10").Value = 51.85751846799701
$ws.Range("R10").Value = 466.7176662119731
$ws.Range("S10").Value = 0.006024770757139146
$ws.Range("T10").Value = 0.006024770757139146
$ws.Range("G11").Value = 15.69885766666667
$ws.Range("H11").Value = 47.096573
$ws.Range("I11").Value = 0.1791719719949428
$ws.Range("J11").Value = 0.1791719719949428
$ws.Range("M11").Value = 37.82684066666667
$ws.Range("O11").Value = 0.3850586149964086
$ws.Range("P11").Value = 0.3850586149964086
$ws.Range("Q11").Value = 593.8381876056785
$ws.Range("R11").Value = 5344.543688451106
$ws.Range("S11").Value = 0.06899171138254799
$ws.Range("T11").Value = 0.06899171138254799
$ws.Range("G12").Value = 15.69885766666667
$ws.Range("H12").Value = 47.096573
$ws.Range("I12").Value = 0.1791719719949428
$ws.Range("J12").Value = 0.1791719719949428
$ws.Range("M12").Value = 9.149395999999999
$ws.Range("N12").Value = 27.448188
$ws.Range("O12").Value = 0.09313634682999644
$ws.Range("P12").Value = 0.09313634682999644
$ws.Range("Q12").Value = 143.6350655399693
$ws.Range("R12").Value = 1292.715589859724
$ws.Range("S12").Value = 0.0166874229259354
$ws.Range("T12").Value = 0.0166874229259354
$ws.Range("G13").Value = 15.69885766666667
$ws.Range("H13").Value = 47.096573
$ws.Range("I13").Value = 0.1791719719949428
$ws.Range("J13").Value = 0.1791719719949428
$ws.Range("M13").Value = 47.95707433333333
$ws.Range("N13").Value = 143.871223
$ws.Range("O13").Value = 0.4881794063849957
$ws.Range("P13").Value = 0.4881794063849957
$ws.Range("Q13").Value = 752.8712840687532
$ws.Range("R13").Value = 6775.841556618779
$ws.Range("S13").Value = 0.08746806692932026
$ws.Range("T13").Value = 0.08746806692932026
$ws.Range("G14").Value = 23.91632366666667
$ws.Range("H14").Value = 71.748971
$ws.Range("I14").Value = 0.2729583874962189
$ws.Range("J14").Value = 0.2729583874962189
$ws.Range("M14").Value = 3.303267
$ws.Range("N14").Value = 9.909801000000002
$ws.Range("O14").Value = 0.03362563178859915
$ws.Range("P14").Value = 0.03362563178859915
$ws.Range("Q14").Value = 79.002002729419
$ws.Range("R14").Value = 711.0180245647711
$ws.Range("S14").Value = 0.009178398231557625
$ws.Range("T14").Value = 0.009178398231557625
$ws.Range("G15").Value = 23.91632366666667
$ws.Range("H15").Value = 71.748971
$ws.Range("I15").Value = 0.2729583874962189
$ws.Range("J15").Value = 0.2729583874962189
$ws.Range("M15").Value = 37.82684066666667
$ws.Range("O15").Value = 0.3850586149964086
$ws.Range("P15").Value = 0.3850586149964086
$ws.Range("Q15").Value = 904.6789646714292
$ws.Range("R15").Value = 8142.110682042862
$ws.Range("S15").Value = 0.1051049786409471
$ws.Range("T15").Value = 0.1051049786409471
$ws.Range("G16").Value = 23.91632366666667
$ws.Range("H16").Value = 71.748971
$ws.Range("I16").Value = 0.2729583874962189
$ws.Range("J16").Value = 0.2729583874962189
$ws.Range("M16").Value = 9.149395999999999
$ws.Range("N16").Value = 27.448188
$ws.Range("O16").Value = 0.09313634682999644
$ws.Range("P16").Value = 0.09313634682999644
$ws.Range("Q16").Value = 218.8199160905053
$ws.Range("R16").Value = 1969.379244814548
$ws.Range("S16").Value = 0.02542234704800441
$ws.Range("T16").Value = 0.02542234704800441
$ws.Range("G17").Value = 23.91632366666667
$ws.Range("H17").Value = 71.748971
$ws.Range("I17").Value = 0.2729583874962189
$ws.Range("J17").Value = 0.2729583874962189
$ws.Range("M17").Value = 47.95707433333333
$ws.Range("N17").Value = 143.871223
$ws.Range("O17").Value = 0.4881794063849957
$ws.Range("P17").Value = 0.4881794063849957
$ws.Range("Q17").Value = 1146.956911862392
$ws.Range("R17").Value = 10322.61220676153
$ws.Range("S17").Value = 0.1332526635757098
$ws.Range("T17").Value = 0.1332526635757098
